# pid4cat_model schema update
# - add slot: pid_schema_version (also: record_version) on PID4CatRecord, and
#   re-order so resource_info/related_identifiers/change_log move after them
# - rename "agent" -> "has_agent" on PID4CatRelation and LogRecord sheets
# - fix: ResourceInfo.resource_category enum CATALYST -> SAMPLE (single-valued rename)
# - add IS_COLLECTED_BY / COLLECTS to the relation_type enum
# - add a new Container sheet with a single contains_pids slot

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. PID4CatRecord: insert record_version + pid_schema_version columns before
#    dc_rights, which shifts resource_info/related_identifiers/change_log right.
# ---------------------------------------------------------------------------
$wsRecord = $wb.Worksheets.Item("PID4CatRecord")

# Current header layout (A:H):
#   A id | B landing_page_url | C status | D resource_info |
#   E related_identifiers | F dc_rights | G curation_contact | H change_log
# Target header layout (A:J):
#   A id | B landing_page_url | C status | D record_version | E pid_schema_version |
#   F dc_rights | G curation_contact | H resource_info | I related_identifiers | J change_log

# Shift the tail (resource_info, related_identifiers, change_log) right by two
# columns first so nothing is clobbered, then fill in the two new columns.
# (NB: read back via Value2 -- Value has been observed to return a stale/
# descriptor string instead of the real cell content on this host.)
$wsRecord.Cells.Item(1, 10).Value = $wsRecord.Cells.Item(1, 8).Value2  # J1 = change_log (was H1)
$wsRecord.Cells.Item(1, 9).Value  = $wsRecord.Cells.Item(1, 5).Value2  # I1 = related_identifiers (was E1)
$wsRecord.Cells.Item(1, 8).Value  = $wsRecord.Cells.Item(1, 4).Value2  # H1 = resource_info (was D1)

$wsRecord.Cells.Item(1, 4).Value = "record_version"
$wsRecord.Cells.Item(1, 5).Value = "pid_schema_version"
# F1 (dc_rights) and G1 (curation_contact) stay put.

# ---------------------------------------------------------------------------
# 2. PID4CatRelation: header rename, and add IS_COLLECTED_BY/COLLECTS to the
#    relation_type dropdown list.
# ---------------------------------------------------------------------------
$wsRelation = $wb.Worksheets.Item("PID4CatRelation")
$wsRelation.Cells.Item(1, 4).Value = "has_agent"

$relationList = '"IS_CITED_BY,CITES,IS_SUPPLEMENT_TO,IS_SUPPLEMENTED_BY,IS_CONTINUED_BY,CONTINUES,HAS_METADATA,IS_METADATA_FOR,HAS_VERSION,IS_VERSION_OF,IS_NEW_VERSION_OF,IS_PREVIOUS_VERSION_OF,IS_PART_OF,HAS_PART,IS_DESCRIBED_BY,DESCRIBES,IS_PUBLISHED_IN,IS_REFERENCED_BY,REFERENCES,IS_DOCUMENTED_BY,DOCUMENTS,IS_COMPILED_BY,COMPILES,IS_VARIANT_FORM_OF,IS_ORIGINAL_FORM_OF,IS_IDENTICAL_TO,IS_DERIVED_FROM,IS_SOURCE_OF,IS_COLLECTED_BY,COLLECTS,IS_REQUIRED_BY,REQUIRES,IS_OBSOLETED_BY,OBSOLETES"'
$relationRange = $wsRelation.Range("A2:A1048576")
$relationRange.Validation.Delete()
$relationRange.Validation.Add(3, 1, 1, $relationList)

# ---------------------------------------------------------------------------
# 3. ResourceInfo: resource_category enum CATALYST -> SAMPLE.
# ---------------------------------------------------------------------------
$wsResourceInfo = $wb.Worksheets.Item("ResourceInfo")
$resourceCategoryList = '"COLLECTION,SAMPLE,MATERIAL,DEVICE,DATAOBJECT"'
$resourceCategoryRange = $wsResourceInfo.Range("C2:C1048576")
$resourceCategoryRange.Validation.Delete()
$resourceCategoryRange.Validation.Add(3, 1, 1, $resourceCategoryList)

# ---------------------------------------------------------------------------
# 4. LogRecord: header rename agent -> has_agent.
# ---------------------------------------------------------------------------
$wsLogRecord = $wb.Worksheets.Item("LogRecord")
$wsLogRecord.Cells.Item(1, 2).Value = "has_agent"

# ---------------------------------------------------------------------------
# 5. Add new "Container" sheet at the end with a single contains_pids slot.
# ---------------------------------------------------------------------------
$wsContainer = $wb.Worksheets.Add()
$wsContainer.Name = "Container"
$wsContainer.Range("A1").Value = "contains_pids"

# Move Container to be the last sheet (right after Agent).
$wsAgent = $wb.Worksheets.Item("Agent")
$wsContainer.Move($null, $wsAgent)
